$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Row=8; Value='Обращение рассмотрено'}
    @{Row=11; Value='Обращение рассмотрено'}
    @{Row=15; Value='Обращение рассмотрено'}
    @{Row=21; Value='Обращение рассмотрено'}
    @{Row=32; Value='Запрос направлен'}
    @{Row=33; Value='Запрос направлен'}
    @{Row=49; Value='Обращение рассмотрено'}
    @{Row=56; Value='Обращение рассмотрено'}
    @{Row=65; Value='Обращение рассмотрено'}
    @{Row=70; Value='Обращение рассмотрено'}
    @{Row=72; Value='Обращение рассмотрено'}
    @{Row=76; Value='Взыскание обращено'}
    @{Row=78; Value='Обращение рассмотрено'}
    @{Row=83; Value='Обращение рассмотрено'}
    @{Row=87; Value='Запрос направлен'}
    @{Row=89; Value='Запрос направлен'}
    @{Row=102; Value='Обращение рассмотрено'}
    @{Row=107; Value='Взыскание обращено'}
    @{Row=111; Value='Обращение рассмотрено'}
    @{Row=114; Value='Запрос направлен'}
    @{Row=115; Value='Постановление вынесено'}
    @{Row=116; Value='Взыскание обращено'}
    @{Row=120; Value='Обращение рассмотрено'}
    @{Row=127; Value='Обращение рассмотрено'}
    @{Row=130; Value='Запрос направлен'}
    @{Row=135; Value='Запрос направлен'}
    @{Row=149; Value='Взыскание обращено'}
    @{Row=151; Value='Обращение рассмотрено'}
    @{Row=158; Value='Обращение рассмотрено'}
    @{Row=174; Value='Обращение рассмотрено'}
    @{Row=176; Value='Взыскание обращено'}
    @{Row=184; Value='Запрос направлен'}
    @{Row=187; Value='Обращение рассмотрено'}
    @{Row=193; Value='Обращение рассмотрено'}
    @{Row=198; Value='Обращение рассмотрено'}
    @{Row=199; Value='Обращение рассмотрено'}
    @{Row=204; Value='Обращение рассмотрено'}
    @{Row=207; Value='Обращение рассмотрено'}
    @{Row=212; Value='Обращение рассмотрено'}
    @{Row=213; Value='Обращение рассмотрено'}
    @{Row=215; Value='Обращение рассмотрено'}
    @{Row=219; Value='Обращение рассмотрено'}
    @{Row=227; Value='Запрос направлен'}
    @{Row=238; Value='Обращение рассмотрено'}
    @{Row=246; Value='Обращение рассмотрено'}
    @{Row=247; Value='Обращение рассмотрено'}
    @{Row=259; Value='Запрос направлен'}
    @{Row=272; Value='Обращение рассмотрено'}
    @{Row=274; Value='Взыскание обращено'}
    @{Row=278; Value='Запрос направлен'}
    @{Row=279; Value='Запрос направлен'}
    @{Row=310; Value='Обращение рассмотрено'}
    @{Row=311; Value='Обращение рассмотрено'}
    @{Row=330; Value='Обращение рассмотрено'}
    @{Row=341; Value='Обращение рассмотрено'}
    @{Row=344; Value='Обращение рассмотрено'}
    @{Row=364; Value='Взыскание обращено'}
    @{Row=397; Value='Обращение рассмотрено'}
    @{Row=398; Value='Обращение рассмотрено'}
    @{Row=399; Value='Обращение рассмотрено'}
    @{Row=411; Value='Взыскание обращено'}
    @{Row=420; Value='Обращение рассмотрено'}
    @{Row=438; Value='Запрос направлен'}
    @{Row=443; Value='Обращение рассмотрено'}
    @{Row=446; Value='Взыскание обращено'}
    @{Row=454; Value='Обращение рассмотрено'}
    @{Row=457; Value='Запрос направлен'}
    @{Row=458; Value='Взыскание обращено'}
    @{Row=462; Value='Взыскание обращено'}
    @{Row=463; Value='Взыскание обращено'}
    @{Row=473; Value='Взыскание обращено'}
    @{Row=507; Value='Обращение рассмотрено'}
    @{Row=508; Value='Обращение рассмотрено'}
    @{Row=512; Value='Обращение рассмотрено'}
    @{Row=528; Value='Обращение рассмотрено'}
    @{Row=529; Value='Взыскание обращено'}
    @{Row=530; Value='Обращение рассмотрено'}
    @{Row=538; Value='Запрос направлен'}
    @{Row=540; Value='Обращение рассмотрено'}
    @{Row=543; Value='Обращение рассмотрено'}
    @{Row=548; Value='Запрос направлен'}
    @{Row=560; Value='Запрос направлен'}
    @{Row=567; Value='Взыскание обращено'}
    @{Row=584; Value='Запрос направлен'}
    @{Row=585; Value='Обращение рассмотрено'}
    @{Row=596; Value='Обращение рассмотрено'}
    @{Row=607; Value='Постановление вынесено'}
    @{Row=620; Value='Обращение рассмотрено'}
    @{Row=647; Value='Обращение рассмотрено'}
    @{Row=649; Value='Запрос направлен'}
    @{Row=691; Value='Обращение рассмотрено'}
    @{Row=692; Value='Обращение рассмотрено'}
    @{Row=698; Value='Запрос направлен'}
    @{Row=699; Value='Обращение рассмотрено'}
    @{Row=710; Value='Обращение рассмотрено'}
    @{Row=719; Value='Обращение рассмотрено'}
    @{Row=733; Value='Обращение рассмотрено'}
    @{Row=748; Value='Запрос направлен'}
    @{Row=757; Value='Взыскание обращено'}
    @{Row=760; Value='Обращение рассмотрено'}
    @{Row=773; Value='Обращение рассмотрено'}
    @{Row=774; Value='Обращение рассмотрено'}
    @{Row=775; Value='Обращение рассмотрено'}
    @{Row=789; Value='Обращение рассмотрено'}
    @{Row=795; Value='Взыскание обращено'}
    @{Row=798; Value='Взыскание обращено'}
    @{Row=799; Value='Взыскание обращено'}
    @{Row=802; Value='Обращение рассмотрено'}
    @{Row=811; Value='Обращение рассмотрено'}
    @{Row=819; Value='Обращение рассмотрено'}
    @{Row=821; Value='Запрос направлен'}
    @{Row=828; Value='Обращение рассмотрено'}
    @{Row=844; Value='Обращение рассмотрено'}
    @{Row=846; Value='Обращение рассмотрено'}
    @{Row=848; Value='Запрос направлен'}
    @{Row=870; Value='Обращение рассмотрено'}
    @{Row=871; Value='Обращение рассмотрено'}
    @{Row=877; Value='Обращение рассмотрено'}
    @{Row=891; Value='Обращение рассмотрено'}
    @{Row=922; Value='Обращение рассмотрено'}
    @{Row=951; Value='Запрос направлен'}
    @{Row=975; Value='Запрос направлен'}
    @{Row=978; Value='Обращение рассмотрено'}
    @{Row=987; Value='Обращение рассмотрено'}
    @{Row=992; Value='Обращение рассмотрено'}
    @{Row=997; Value='Обращение рассмотрено'}
)

foreach ($change in $changes) {
    $ws.Cells.Item($change.Row, 3).Value = $change.Value
}
